$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 380 entirely; Excel shifts all rows below it up by one,
# matching the commit's removal of the "「黙ってもらえます？」" post row.
$ws.Rows.Item(380).Delete()
